$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 88.8032398676725
$ws.Cells.Item(2, 5).Value = 88.323974609375
$ws.Cells.Item(2, 6).Value = 96.38156510813448
$ws.Cells.Item(2, 7).Value = 87.79479080641208
$ws.Cells.Item(2, 8).Value = 24003656
$ws.Cells.Item(2, 9).Value = "FICO"

$ws.Cells.Item(3, 4).Value = 91.50983758235208
$ws.Cells.Item(3, 5).Value = 90.57107543945312
$ws.Cells.Item(3, 6).Value = 97.43206153799773
$ws.Cells.Item(3, 7).Value = 89.49248987674569
$ws.Cells.Item(3, 8).Value = 24003656
$ws.Cells.Item(3, 9).Value = "FICO"

$ws.Cells.Item(4, 4).Value = 84.41925966923348
$ws.Cells.Item(4, 5).Value = 92.27082061767578
$ws.Cells.Item(4, 6).Value = 94.6482623532052
$ws.Cells.Item(4, 7).Value = 82.20164217095861
$ws.Cells.Item(4, 8).Value = 24003656
$ws.Cells.Item(4, 9).Value = "FICO"

$ws.Cells.Item(5, 4).Value = 93.24943876866033
$ws.Cells.Item(5, 5).Value = 95.4875030517578
$ws.Cells.Item(5, 6).Value = 98.87457617196462
$ws.Cells.Item(5, 7).Value = 80.13076779115008
$ws.Cells.Item(5, 8).Value = 24003656
$ws.Cells.Item(5, 9).Value = "FICO"

$ws.Cells.Item(6, 4).Value = 105.0200299744409
$ws.Cells.Item(6, 5).Value = 106.6389541625977
$ws.Cells.Item(6, 6).Value = 109.1772640072961
$ws.Cells.Item(6, 7).Value = 102.7015748789839
$ws.Cells.Item(6, 8).Value = 24003656
$ws.Cells.Item(6, 9).Value = "FICO"

$ws.Cells.Item(7, 4).Value = 112.6859743556066
$ws.Cells.Item(7, 5).Value = 126.5793151855469
$ws.Cells.Item(7, 6).Value = 130.0976327187569
$ws.Cells.Item(7, 7).Value = 111.6764638247111
$ws.Cells.Item(7, 8).Value = 24003656
$ws.Cells.Item(7, 9).Value = "FICO"

$ws.Cells.Item(8, 4).Value = 123.6999527063858
$ws.Cells.Item(8, 5).Value = 120.6409454345703
$ws.Cells.Item(8, 6).Value = 125.9592234530672
$ws.Cells.Item(8, 7).Value = 118.9414986451007
$ws.Cells.Item(8, 8).Value = 24003656
$ws.Cells.Item(8, 9).Value = "FICO"

$ws.Cells.Item(9, 4).Value = 119.8015757470624
$ws.Cells.Item(9, 5).Value = 123.2810440063477
$ws.Cells.Item(9, 6).Value = 126.6805151077219
$ws.Cells.Item(9, 7).Value = 118.9317067751857
$ws.Cells.Item(9, 8).Value = 24003656
$ws.Cells.Item(9, 9).Value = "FICO"

$ws.Cells.Item(10, 4).Value = 128.5
$ws.Cells.Item(10, 5).Value = 135.4799957275391
$ws.Cells.Item(10, 6).Value = 137.5700073242188
$ws.Cells.Item(10, 7).Value = 125.7099990844727
$ws.Cells.Item(10, 8).Value = 24003656
$ws.Cells.Item(10, 9).Value = "FICO"

$ws.Cells.Item(11, 4).Value = 139.8899993896484
$ws.Cells.Item(11, 5).Value = 142.5500030517578
$ws.Cells.Item(11, 6).Value = 146.3399963378906
$ws.Cells.Item(11, 7).Value = 136.6000061035156
$ws.Cells.Item(11, 8).Value = 24003656
$ws.Cells.Item(11, 9).Value = "FICO"

$ws.Cells.Item(12, 4).Value = 140.6199951171875
$ws.Cells.Item(12, 5).Value = 145.1600036621094
$ws.Cells.Item(12, 6).Value = 149
$ws.Cells.Item(12, 7).Value = 140.6199951171875
$ws.Cells.Item(12, 8).Value = 24003656
$ws.Cells.Item(12, 9).Value = "FICO"

$ws.Cells.Item(13, 4).Value = 153.4100036621094
$ws.Cells.Item(13, 5).Value = 172.6600036621094
$ws.Cells.Item(13, 6).Value = 177.9100036621094
$ws.Cells.Item(13, 7).Value = 152.4700012207031
$ws.Cells.Item(13, 8).Value = 24003656
$ws.Cells.Item(13, 9).Value = "FICO"

$ws.Cells.Item(14, 4).Value = 168.4400024414062
$ws.Cells.Item(14, 5).Value = 173.1799926757812
$ws.Cells.Item(14, 6).Value = 177.9199981689453
$ws.Cells.Item(14, 7).Value = 162.5800018310547
$ws.Cells.Item(14, 8).Value = 24003656
$ws.Cells.Item(14, 9).Value = "FICO"

$ws.Cells.Item(15, 4).Value = 193.1100006103516
$ws.Cells.Item(15, 5).Value = 201.4600067138672
$ws.Cells.Item(15, 6).Value = 212.4400024414062
$ws.Cells.Item(15, 7).Value = 193.1100006103516
$ws.Cells.Item(15, 8).Value = 24003656
$ws.Cells.Item(15, 9).Value = "FICO"

$ws.Cells.Item(16, 4).Value = 229.9400024414062
$ws.Cells.Item(16, 5).Value = 192.7100067138672
$ws.Cells.Item(16, 6).Value = 231.9600067138672
$ws.Cells.Item(16, 7).Value = 185.229995727539
$ws.Cells.Item(16, 8).Value = 24003656
$ws.Cells.Item(16, 9).Value = "FICO"

$ws.Cells.Item(17, 4).Value = 184.1600036621093
$ws.Cells.Item(17, 5).Value = 225.1999969482422
$ws.Cells.Item(17, 6).Value = 227.1100006103516
$ws.Cells.Item(17, 7).Value = 178.4199981689453
$ws.Cells.Item(17, 8).Value = 24003656
$ws.Cells.Item(17, 9).Value = "FICO"

$ws.Cells.Item(18, 4).Value = 273.7799987792969
$ws.Cells.Item(18, 5).Value = 279.75
$ws.Cells.Item(18, 6).Value = 287.2900085449219
$ws.Cells.Item(18, 7).Value = 269.2000122070312
$ws.Cells.Item(18, 8).Value = 24003656
$ws.Cells.Item(18, 9).Value = "FICO"

$ws.Cells.Item(19, 4).Value = 318.7000122070312
$ws.Cells.Item(19, 5).Value = 347.4200134277344
$ws.Cells.Item(19, 6).Value = 354.4200134277344
$ws.Cells.Item(19, 7).Value = 310.010009765625
$ws.Cells.Item(19, 8).Value = 24003656
$ws.Cells.Item(19, 9).Value = "FICO"

$ws.Cells.Item(20, 4).Value = 303.25
$ws.Cells.Item(20, 5).Value = 304.0400085449219
$ws.Cells.Item(20, 6).Value = 315.489990234375
$ws.Cells.Item(20, 7).Value = 281.5599975585937
$ws.Cells.Item(20, 8).Value = 24003656
$ws.Cells.Item(20, 9).Value = "FICO"

$ws.Cells.Item(21, 4).Value = 377
$ws.Cells.Item(21, 5).Value = 402.3800048828125
$ws.Cells.Item(21, 6).Value = 420.2000122070313
$ws.Cells.Item(21, 7).Value = 375.25
$ws.Cells.Item(21, 8).Value = 24003656
$ws.Cells.Item(21, 9).Value = "FICO"

$ws.Cells.Item(22, 4).Value = 296.3999938964844
$ws.Cells.Item(22, 5).Value = 352.9400024414062
$ws.Cells.Item(22, 6).Value = 364.260009765625
$ws.Cells.Item(22, 7).Value = 259.3699951171875
$ws.Cells.Item(22, 8).Value = 24003656
$ws.Cells.Item(22, 9).Value = "FICO"

$ws.Cells.Item(23, 4).Value = 416.5199890136719
$ws.Cells.Item(23, 5).Value = 439.1900024414063
$ws.Cells.Item(23, 6).Value = 442.6499938964844
$ws.Cells.Item(23, 7).Value = 388.0199890136719
$ws.Cells.Item(23, 8).Value = 24003656
$ws.Cells.Item(23, 9).Value = "FICO"

$ws.Cells.Item(24, 4).Value = 431.6199951171875
$ws.Cells.Item(24, 5).Value = 391.4500122070313
$ws.Cells.Item(24, 6).Value = 450.25
$ws.Cells.Item(24, 7).Value = 380
$ws.Cells.Item(24, 8).Value = 24003656
$ws.Cells.Item(24, 9).Value = "FICO"

$ws.Cells.Item(25, 4).Value = 511.989990234375
$ws.Cells.Item(25, 5).Value = 450.1099853515625
$ws.Cells.Item(25, 6).Value = 514.75
$ws.Cells.Item(25, 7).Value = 436.6600036621094
$ws.Cells.Item(25, 8).Value = 24003656
$ws.Cells.Item(25, 9).Value = "FICO"

$ws.Cells.Item(26, 4).Value = 492.260009765625
$ws.Cells.Item(26, 5).Value = 521.4099731445312
$ws.Cells.Item(26, 6).Value = 547.5700073242188
$ws.Cells.Item(26, 7).Value = 490
$ws.Cells.Item(26, 8).Value = 24003656
$ws.Cells.Item(26, 9).Value = "FICO"

$ws.Cells.Item(27, 4).Value = 503.5299987792969
$ws.Cells.Item(27, 5).Value = 523.9099731445312
$ws.Cells.Item(27, 6).Value = 553.969970703125
$ws.Cells.Item(27, 7).Value = 500.3200073242188
$ws.Cells.Item(27, 8).Value = 24003656
$ws.Cells.Item(27, 9).Value = "FICO"

$ws.Cells.Item(28, 4).Value = 400.4100036621094
$ws.Cells.Item(28, 5).Value = 398.2000122070313
$ws.Cells.Item(28, 6).Value = 419.9500122070313
$ws.Cells.Item(28, 7).Value = 392.0599975585938
$ws.Cells.Item(28, 8).Value = 24003656
$ws.Cells.Item(28, 9).Value = "FICO"

$ws.Cells.Item(29, 4).Value = 437.1499938964844
$ws.Cells.Item(29, 5).Value = 494.989990234375
$ws.Cells.Item(29, 6).Value = 499.0400085449219
$ws.Cells.Item(29, 7).Value = 413.9200134277344
$ws.Cells.Item(29, 8).Value = 24003656
$ws.Cells.Item(29, 9).Value = "FICO"

$ws.Cells.Item(30, 4).Value = 466.0499877929688
$ws.Cells.Item(30, 5).Value = 373.510009765625
$ws.Cells.Item(30, 6).Value = 468.6700134277344
$ws.Cells.Item(30, 7).Value = 367.4500122070313
$ws.Cells.Item(30, 8).Value = 24003656
$ws.Cells.Item(30, 9).Value = "FICO"

$ws.Cells.Item(31, 4).Value = 398.3299865722656
$ws.Cells.Item(31, 5).Value = 462.0299987792969
$ws.Cells.Item(31, 6).Value = 471.9400024414063
$ws.Cells.Item(31, 7).Value = 397.4100036621094
$ws.Cells.Item(31, 8).Value = 24003656
$ws.Cells.Item(31, 9).Value = "FICO"

$ws.Cells.Item(32, 4).Value = 415.1600036621094
$ws.Cells.Item(32, 5).Value = 478.8399963378906
$ws.Cells.Item(32, 6).Value = 487.8299865722656
$ws.Cells.Item(32, 7).Value = 389.8399963378906
$ws.Cells.Item(32, 8).Value = 24003656
$ws.Cells.Item(32, 9).Value = "FICO"

$ws.Cells.Item(33, 4).Value = 608.6500244140625
$ws.Cells.Item(33, 5).Value = 665.9500122070312
$ws.Cells.Item(33, 6).Value = 676.97998046875
$ws.Cells.Item(33, 7).Value = 575.3900146484375
$ws.Cells.Item(33, 8).Value = 24003656
$ws.Cells.Item(33, 9).Value = "FICO"

$ws.Cells.Item(34, 4).Value = 700
$ws.Cells.Item(34, 5).Value = 727.9500122070312
$ws.Cells.Item(34, 6).Value = 758
$ws.Cells.Item(34, 7).Value = 664.4099731445312
$ws.Cells.Item(34, 8).Value = 24003656
$ws.Cells.Item(34, 9).Value = "FICO"

$ws.Cells.Item(35, 4).Value = 805.9000244140625
$ws.Cells.Item(35, 5).Value = 837.969970703125
$ws.Cells.Item(35, 6).Value = 860
$ws.Cells.Item(35, 7).Value = 776.6300048828125
$ws.Cells.Item(35, 8).Value = 24003656
$ws.Cells.Item(35, 9).Value = "FICO"

$ws.Cells.Item(36, 4).Value = 865
$ws.Cells.Item(36, 5).Value = 845.8699951171875
$ws.Cells.Item(36, 6).Value = 940.0999755859376
$ws.Cells.Item(36, 7).Value = 811.989990234375
$ws.Cells.Item(36, 8).Value = 24003656
$ws.Cells.Item(36, 9).Value = "FICO"

$ws.Cells.Item(38, 4).Value = 1249.609985351562
$ws.Cells.Item(38, 5).Value = 1133.329956054688
$ws.Cells.Item(38, 6).Value = 1266.859985351562
$ws.Cells.Item(38, 7).Value = 1105.650024414062
$ws.Cells.Item(38, 8).Value = 24003656
$ws.Cells.Item(38, 9).Value = "FICO"

$ws.Cells.Item(39, 4).Value = 1495.199951171875
$ws.Cells.Item(39, 5).Value = 1600
$ws.Cells.Item(39, 6).Value = 1658.030029296875
$ws.Cells.Item(39, 7).Value = 1454.410034179688
$ws.Cells.Item(39, 8).Value = 24003656
$ws.Cells.Item(39, 9).Value = "FICO"

$ws.Cells.Item(40, 4).Value = 1953.569946289062
$ws.Cells.Item(40, 5).Value = 1993.109985351562
$ws.Cells.Item(40, 6).Value = 2103.699951171875
$ws.Cells.Item(40, 7).Value = 1886.199951171875
$ws.Cells.Item(40, 8).Value = 24003656
$ws.Cells.Item(40, 9).Value = "FICO"

$ws.Cells.Item(41, 4).Value = 2005
$ws.Cells.Item(41, 5).Value = 1873.56005859375
$ws.Cells.Item(41, 6).Value = 2066.679931640625
$ws.Cells.Item(41, 7).Value = 1787.569946289062
$ws.Cells.Item(41, 8).Value = 24003656
$ws.Cells.Item(41, 9).Value = "FICO"

$ws.Cells.Item(42, 4).Value = 1835.099975585937
$ws.Cells.Item(42, 5).Value = 1989.680053710937
$ws.Cells.Item(42, 6).Value = 2016.780029296875
$ws.Cells.Item(42, 7).Value = 1585.180053710938
$ws.Cells.Item(42, 8).Value = 24003656
$ws.Cells.Item(42, 9).Value = "FICO"

$ws.Cells.Item(43, 4).Value = 1806.420043945312
$ws.Cells.Item(43, 5).Value = 1436.719970703125
$ws.Cells.Item(43, 6).Value = 1905.329956054688
$ws.Cells.Item(43, 7).Value = 1354.47998046875
$ws.Cells.Item(43, 8).Value = 24003656
$ws.Cells.Item(43, 9).Value = "FICO"
